# Centralize Services.java / update core objects:
# Append two new rows of test case data to the TestCaseMaster sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseMaster")

# New row 10 - JDServices CRUD test case
$ws.Cells.Item(10, 1).Value = "9"
$ws.Cells.Item(10, 2).Value = "JDServices"
$ws.Cells.Item(10, 3).Value = "CRUD"
$ws.Cells.Item(10, 4).Value = "JDservices.CRUD"
$ws.Cells.Item(10, 5).Value = "ServiceTests//DummySheet.xlsx"
$ws.Cells.Item(10, 6).Value = "create"
$ws.Cells.Item(10, 7).Value = "All"

# New row 11 - LocationsServices CRUD test case
$ws.Cells.Item(11, 1).Value = "10"
$ws.Cells.Item(11, 4).Value = "LocationsServices.CRUD"
$ws.Cells.Item(11, 2).Value = "LocationsServices"
$ws.Cells.Item(11, 3).Value = "CRUD"
$ws.Cells.Item(11, 5).Value = "ServiceTests//DummySheet.xlsx"
$ws.Cells.Item(11, 6).Value = "create"
$ws.Cells.Item(11, 7).Value = "All"

# Match the styling already used on columns A and G (text number format)
$ws.Range("A10:A11").NumberFormat = "@"
$ws.Range("G10:G11").NumberFormat = "@"

# Update the active selection to mirror the recorded cursor position
$ws.Range("H11").Select()
